$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1030.6364
$ws.Range("I19").Value = 1329
$ws.Range("J19").Value = 782
$ws.Range("K19").Value = 1329
$ws.Range("L19").Value = 782
$ws.Range("M19").Value = -1154
$ws.Range("N19").Value = -1132

$ws.Range("H74").Value = 1989
$ws.Range("I74").Value = 1989
$ws.Range("K74").Value = 1989
$ws.Range("M74").Value = -1053

$ws.Range("H77").Value = 1989
$ws.Range("I77").Value = 1989
$ws.Range("K77").Value = 9945
$ws.Range("M77").Value = -5265

$ws.Range("H100").Value = 1622.1818
$ws.Range("J100").Value = 997.5
$ws.Range("L100").Value = 997.5
$ws.Range("N100").Value = -2079.5

$ws.Range("H107").Value = 1832.6666
$ws.Range("I107").Value = 1729.2
$ws.Range("K107").Value = 1729.2
$ws.Range("M107").Value = 190.8

$ws.Range("H113").Value = 4449.5835
$ws.Range("J113").Value = 3438.3333
$ws.Range("L113").Value = 3438.3333
$ws.Range("N113").Value = -9946.3333

$ws.Range("H137").Value = 1380.75
$ws.Range("I137").Value = 1063.1666
$ws.Range("K137").Value = 3189.4998
$ws.Range("M137").Value = -639.4998000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5011
$ws.Range("I61").Value = 5011
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5011
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4799
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

$ws.Range("H97").Value = 737.3333
$ws.Range("I97").Value = 457.6154
$ws.Range("K97").Value = 457.6154
$ws.Range("M97").Value = 38.38459999999998

$ws.Range("H102").Value = 4090
$ws.Range("I102").Value = 1827.1428
$ws.Range("J102").Value = 8050
$ws.Range("K102").Value = 1827.1428
$ws.Range("L102").Value = 8050
$ws.Range("M102").Value = -205.1428000000001
$ws.Range("N102").Value = -11294

$ws.Range("H122").Value = 999.5
$ws.Range("I122").Value = 999.5
$ws.Range("K122").Value = 2998.5
$ws.Range("M122").Value = -548.5

$ws.Range("H130").Value = 13000
$ws.Range("J130").Value = 13000
$ws.Range("L130").Value = 13000
$ws.Range("N130").Value = -23040

$ws.Range("H136").Value = 5011
$ws.Range("I136").Value = 5011
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 15033
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -12483
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1858.3334
$ws.Range("I22").Value = 1858.3334
$ws.Range("K22").Value = 1858.3334
$ws.Range("M22").Value = -1685.3334

$ws.Range("H107").Value = 2384
$ws.Range("I107").Value = 1622.4
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 1622.4
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = 297.5999999999999
$ws.Range("N107").Value = -13840

$ws.Range("H134").Value = 2990
$ws.Range("I134").Value = 2990
$ws.Range("K134").Value = 8970
$ws.Range("M134").Value = -6435

$ws.Range("H135").Value = 1198296
$ws.Range("J135").Value = 1198296
$ws.Range("L135").Value = 1198296
$ws.Range("N135").Value = -1208436

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4108.59
$ws.Range("I31").Value = 1955.9286
$ws.Range("K31").Value = 1955.9286
$ws.Range("M31").Value = -1660.9286

$ws.Range("H34").Value = 4108.59
$ws.Range("I34").Value = 1955.9286
$ws.Range("K34").Value = 1955.9286
$ws.Range("M34").Value = -1753.9286

$ws.Range("H51").Value = 67951.664
$ws.Range("I51").Value = 15090
$ws.Range("K51").Value = 15090
$ws.Range("M51").Value = -14354

$ws.Range("H61").Value = 67951.664
$ws.Range("I61").Value = 15090
$ws.Range("K61").Value = 15090
$ws.Range("M61").Value = -14742

$ws.Range("H132").Value = 1996.25
$ws.Range("I132").Value = 1996.25
$ws.Range("K132").Value = 5988.75
$ws.Range("M132").Value = -3458.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83.1579
$ws.Range("J2").Value = 58.75
$ws.Range("L2").Value = 352.5
$ws.Range("N2").Value = -578.5

$ws.Range("H10").Value = 74.75
$ws.Range("J10").Value = 150.75
$ws.Range("L10").Value = 452.25
$ws.Range("N10").Value = -730.25

$ws.Range("H12").Value = 111
$ws.Range("J12").Value = 138.5
$ws.Range("L12").Value = 415.5
$ws.Range("N12").Value = -761.5

$ws.Range("H131").Value = 2362.375
$ws.Range("I131").Value = 2399
$ws.Range("J131").Value = 2357.1428
$ws.Range("K131").Value = 7197
$ws.Range("L131").Value = 7071.428400000001
$ws.Range("M131").Value = -2157
$ws.Range("N131").Value = -17151.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 296.41666
$ws.Range("I2").Value = 339.75
$ws.Range("J2").Value = 274.75
$ws.Range("K2").Value = 339.75
$ws.Range("L2").Value = 274.75
$ws.Range("M2").Value = -226.75
$ws.Range("N2").Value = -500.75

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H40").Value = 19600
$ws.Range("J40").Value = 19600
$ws.Range("L40").Value = 19600
$ws.Range("N40").Value = -19902

$ws.Range("H42").Value = 106998.336
$ws.Range("J42").Value = 106998.336
$ws.Range("L42").Value = 106998.336
$ws.Range("N42").Value = -107968.336

$ws.Range("H43").Value = 14766.556
$ws.Range("J43").Value = 19747
$ws.Range("L43").Value = 19747
$ws.Range("N43").Value = -20049

$ws.Range("H57").Value = 5500
$ws.Range("I57").Value = 5500
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 5500
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -4680
$ws.Range("N57").ClearContents()

$ws.Range("H70").Value = 3889
$ws.Range("J70").Value = 3664
$ws.Range("L70").Value = 3664
$ws.Range("N70").Value = -4204

$ws.Range("H73").Value = 3889
$ws.Range("J73").Value = 3664
$ws.Range("L73").Value = 3664
$ws.Range("N73").Value = -5536

$ws.Range("H80").Value = 3606
$ws.Range("I80").Value = 3990
$ws.Range("J80").Value = 3414
$ws.Range("K80").Value = 3990
$ws.Range("L80").Value = 3414
$ws.Range("M80").Value = -2992
$ws.Range("N80").Value = -5410

$ws.Range("H83").Value = 3606
$ws.Range("I83").Value = 3990
$ws.Range("J83").Value = 3414
$ws.Range("K83").Value = 19950
$ws.Range("L83").Value = 17070
$ws.Range("M83").Value = -14958
$ws.Range("N83").Value = -27054

$ws.Range("H102").Value = 2184.625
$ws.Range("I102").Value = 2062.2173
$ws.Range("K102").Value = 2062.2173
$ws.Range("M102").Value = -440.2172999999998

$ws.Range("H115").Value = 106998.336
$ws.Range("J115").Value = 106998.336
$ws.Range("L115").Value = 106998.336
$ws.Range("N115").Value = -109348.336

$ws.Range("H122").Value = 2889.1428
$ws.Range("I122").Value = 2074.5
$ws.Range("K122").Value = 6223.5
$ws.Range("M122").Value = -3773.5

$ws.Range("H132").Value = 2581.8333
$ws.Range("I132").Value = 2581.8333
$ws.Range("K132").Value = 7745.499899999999
$ws.Range("M132").Value = -5215.499899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6538.1816
$ws.Range("I7").Value = 6498.4287
$ws.Range("K7").Value = 6498.4287
$ws.Range("M7").Value = -6386.4287

$ws.Range("H20").Value = 28036.324
$ws.Range("I20").Value = 1096.5625
$ws.Range("K20").Value = 1096.5625
$ws.Range("M20").Value = -870.5625

$ws.Range("H126").Value = 6538.1816
$ws.Range("I126").Value = 6498.4287
$ws.Range("K126").Value = 19495.2861
$ws.Range("M126").Value = -17025.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5400
$ws.Range("I81").Value = 800
$ws.Range("K81").Value = 1600
$ws.Range("M81").Value = -539

$ws.Range("H84").Value = 5400
$ws.Range("I84").Value = 800
$ws.Range("K84").Value = 8000
$ws.Range("M84").Value = -2696

$ws.Range("H100").Value = 1056.8
$ws.Range("I100").Value = 508.5
$ws.Range("K100").Value = 1017
$ws.Range("M100").Value = -476

$ws.Range("H107").Value = 1405.5
$ws.Range("I107").Value = 1967.6666
$ws.Range("J107").Value = 562.25
$ws.Range("K107").Value = 5902.9998
$ws.Range("L107").Value = 1686.75
$ws.Range("M107").Value = -3982.9998
$ws.Range("N107").Value = -5526.75

$ws.Range("H122").Value = 5618.375
$ws.Range("I122").Value = 3629.4
$ws.Range("K122").Value = 10888.2
$ws.Range("M122").Value = -8438.200000000001

